# Apply cryptos list update (prices/volumes refreshed, a few rows re-ranked)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "42.760.56"
$ws.Cells.Item(2, 5).Value = "  -0.45%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.294.46"
$ws.Cells.Item(3, 5).Value = "  -0.21%  "
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "303.90"
$ws.Cells.Item(5, 5).Value = "  +1.46%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "96.55"
$ws.Cells.Item(6, 5).Value = "  -0.96%  "
$ws.Cells.Item(7, 5).Value = "  -2.02%  "
$ws.Cells.Item(8, 5).Value = "  -0.07%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.498"
$ws.Cells.Item(9, 5).Value = "  -1.79%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "35.06"
$ws.Cells.Item(10, 5).Value = "  -2.18%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0783"
$ws.Cells.Item(11, 5).Value = "  -0.68%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "18.77"
$ws.Cells.Item(12, 5).Value = "  +5.90%  "
$ws.Cells.Item(13, 5).Value = "  +2.01%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.84"
$ws.Cells.Item(14, 5).Value = "  +0.90%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "2.651.73"
$ws.Cells.Item(15, 5).Value = "  -0.25%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.308.48"
$ws.Cells.Item(16, 5).Value = "  -0.18%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.775"
$ws.Cells.Item(17, 5).Value = "  -0.34%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "42.680.40"
$ws.Cells.Item(18, 5).Value = "  -0.56%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.84"
$ws.Cells.Item(19, 5).Value = "  +2.01%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0₃0895"
$ws.Cells.Item(20, 5).Value = "  -1.44%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.99"
$ws.Cells.Item(21, 5).Value = "  -1.67%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "67.22"
$ws.Cells.Item(22, 5).Value = "  -1.31%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "236.25"
$ws.Cells.Item(23, 5).Value = "  -2.32%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.15"
$ws.Cells.Item(24, 5).Value = "  +0.32%  "
$ws.Cells.Item(25, 5).Value = "  +0.10%  "
$ws.Cells.Item(26, 5).Value = "  -1.55%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "25.00"
$ws.Cells.Item(27, 5).Value = "  -0.35%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "167.09"
$ws.Cells.Item(28, 5).Value = "  +0.33%  "
$ws.Cells.Item(29, 5).Value = "  +0.91%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.01"
$ws.Cells.Item(30, 5).Value = "  -0.43%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "32.98"
$ws.Cells.Item(31, 5).Value = "  +0.22%  "
$ws.Cells.Item(32, 5).Value = "  +0.04%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "17.89"
$ws.Cells.Item(33, 5).Value = "  +1.33%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.97"
$ws.Cells.Item(34, 5).Value = "  -0.59%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.48"
$ws.Cells.Item(35, 5).Value = "  -4.82%  "
$ws.Cells.Item(36, 5).Value = "  -1.67%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0684"
$ws.Cells.Item(37, 5).Value = "  -0.13%  "
$ws.Cells.Item(38, 5).Value = "  -0.20%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.75"
$ws.Cells.Item(39, 5).Value = "  -0.68%  "
$ws.Cells.Item(40, 5).Value = "  -0.91%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.69"
$ws.Cells.Item(41, 5).Value = "  -2.17%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.992.07"
$ws.Cells.Item(42, 5).Value = "  -0.51%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0280"
$ws.Cells.Item(43, 5).Value = "  -2.37%  "
$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "18.36"
$ws.Cells.Item(44, 5).Value = "  +5.28%  "
$ws.Cells.Item(45, 2).Value = "FraxShare"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "10.20"
$ws.Cells.Item(45, 5).Value = "  +0.38%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.13"
$ws.Cells.Item(46, 5).Value = "  -0.85%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.78"
$ws.Cells.Item(47, 5).Value = "  +0.33%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.89"
$ws.Cells.Item(48, 5).Value = "  -0.28%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "53.68"
$ws.Cells.Item(49, 5).Value = "  +0.53%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.518.46"
$ws.Cells.Item(50, 5).Value = "  -0.30%  "
$ws.Cells.Item(51, 2).Value = "TrustWalletToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.13"
$ws.Cells.Item(51, 5).Value = "  +1.68%  "
